$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update current-balance ratios for the three items
$ws.Range("H7").Value = "1:0"
$ws.Range("H8").Value = "0:-1"
$ws.Range("H9").Value = "0:0"

# Update the generation timestamp
$ws.Range("A11").Value = "Saturday, 24 May, 2025 9:46 AM"
